$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the old "units" row (row 2: (m3/s), (MW), (MW), (GWh), (GWh), (GWh))
# This shifts all data rows (old row 3..99) up by one.
$ws.Rows.Item(2).Delete()

# Overwrite row 1 with the new header labels.
$ws.Cells.Item(1, 1).Value = "idx"
$ws.Cells.Item(1, 2).Value = "idx2"
$ws.Cells.Item(1, 3).Value = "Name"
$ws.Cells.Item(1, 4).Value = "Date Start"
$ws.Cells.Item(1, 5).Value = "Date End"
$ws.Cells.Item(1, 6).Value = "(m3/s)"
$ws.Cells.Item(1, 7).Value = "(MW1)"
$ws.Cells.Item(1, 8).Value = "(MW2)"
$ws.Cells.Item(1, 9).Value = "(GWh) Winter"
$ws.Cells.Item(1, 10).Value = "(GWh) Summer"
$ws.Cells.Item(1, 11).Value = "(GWh) Year"

$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Active selection moves to A2:K2 (the first data row) per the sheetView.
$ws.Range("A2:K2").Select()
